$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 for the "loss" level, continuing the existing table
$ws.Range("B12").Value = "loss"
$ws.Range("G12").Value = 5

# Column H: "1.5x the winnings", column I: "double the winnings" -- row 2
# was entered first as a one-off formula, then H3:H12/I3:I12 were filled
# down as their own block
$ws.Range("H2").Formula = "=G2*1.5"
$ws.Range("I2").Formula = "=G2*2"
$ws.Range("H3:H12").Formula = "=G3*1.5"
$ws.Range("I3:I12").Formula = "=G3*2"

# ... except the rows where the author overtyped the computed value with a
# hand-rounded whole number, breaking the shared-formula chain there
$ws.Range("H7").Value = 68
$ws.Range("H9").Value = 38
$ws.Range("H12").Value = 8

# Leave the cursor where the editing session ended
$ws.Range("I12").Select()
